$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: League changes to "Italian Serie A"; Time/Home/Away updated to the new fixture
$ws.Cells.Item(2, 1).Value = "Italian Serie A"
$ws.Cells.Item(2, 3).Value = "14:30:00"
$ws.Cells.Item(2, 4).Value = "Verona"
$ws.Cells.Item(2, 5).Value = "Bologna"

# Row 3: German Bundesliga, Augsburg vs Union Berlin
$ws.Cells.Item(3, 1).Value = "German Bundesliga"
$ws.Cells.Item(3, 3).Value = "16:30:00"
$ws.Cells.Item(3, 4).Value = "Augsburg"
$ws.Cells.Item(3, 5).Value = "Union Berlin"
$ws.Cells.Item(3, 6).Value = 3.15
$ws.Cells.Item(3, 7).Value = 3.2
$ws.Cells.Item(3, 8).Value = 2.96
$ws.Cells.Item(3, 9).Value = 2.98
$ws.Cells.Item(3, 10).Value = 2.9
$ws.Cells.Item(3, 11).Value = 2.92
$ws.Cells.Item(3, 12).Value = 2.66
$ws.Cells.Item(3, 13).Value = 1.2
$ws.Cells.Item(3, 14).Value = 2.24
$ws.Cells.Item(3, 15).Value = 1.79
$ws.Cells.Item(3, 16).Value = 1.34
$ws.Cells.Item(3, 17).Value = 3.85
$ws.Cells.Item(3, 18).Value = 1.11
$ws.Cells.Item(3, 19).Value = 9.4
$ws.Cells.Item(3, 20).Value = 2.56
$ws.Cells.Item(3, 21).Value = 1.58
$ws.Cells.Item(3, 22).Value = 1.51
$ws.Cells.Item(3, 23).Value = 1.45
$ws.Cells.Item(3, 24).Value = 6.2
$ws.Cells.Item(3, 25).Value = 6.8
$ws.Cells.Item(3, 26).Value = 15.5
$ws.Cells.Item(3, 27).Value = 60
$ws.Cells.Item(3, 28).Value = 7.4
$ws.Cells.Item(3, 29).Value = 7
$ws.Cells.Item(3, 30).Value = 16.5
$ws.Cells.Item(3, 31).Value = 60
$ws.Cells.Item(3, 32).Value = 17.5
$ws.Cells.Item(3, 33).Value = 16.5
$ws.Cells.Item(3, 34).Value = 36
$ws.Cells.Item(3, 35).Value = 150
$ws.Cells.Item(3, 36).Value = 75
$ws.Cells.Item(3, 37).Value = 70
$ws.Cells.Item(3, 38).Value = 160
$ws.Cells.Item(3, 39).Value = 470
$ws.Cells.Item(3, 40).Value = 140
$ws.Cells.Item(3, 41).Value = 110

# Row 4: Italian Serie A, Como vs AC Milan
$ws.Cells.Item(4, 1).Value = "Italian Serie A"
$ws.Cells.Item(4, 3).Value = "16:45:00"
$ws.Cells.Item(4, 4).Value = "Como"
$ws.Cells.Item(4, 5).Value = "AC Milan"
$ws.Cells.Item(4, 6).Value = 1.74
$ws.Cells.Item(4, 7).Value = 1.75
$ws.Cells.Item(4, 8).Value = 5.2
$ws.Cells.Item(4, 9).Value = 5.3
$ws.Cells.Item(4, 10).Value = 4.2
$ws.Cells.Item(4, 11).Value = 4.3
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 0
$ws.Cells.Item(4, 14).Value = 11
$ws.Cells.Item(4, 15).Value = 1.09
$ws.Cells.Item(4, 16).Value = 3.2
$ws.Cells.Item(4, 17).Value = 1.44
$ws.Cells.Item(4, 18).Value = 1.71
$ws.Cells.Item(4, 19).Value = 2.36
$ws.Cells.Item(4, 20).Value = 1.35
$ws.Cells.Item(4, 21).Value = 3.7
$ws.Cells.Item(4, 22).Value = 1.23
$ws.Cells.Item(4, 23).Value = 2.34
$ws.Cells.Item(4, 24).Value = 1000
$ws.Cells.Item(4, 25).Value = 1000
$ws.Cells.Item(4, 26).Value = 1000
$ws.Cells.Item(4, 27).Value = 1000
$ws.Cells.Item(4, 28).Value = 12
$ws.Cells.Item(4, 29).Value = 8.2
$ws.Cells.Item(4, 30).Value = 12.5
$ws.Cells.Item(4, 31).Value = 36
$ws.Cells.Item(4, 32).Value = 10.5
$ws.Cells.Item(4, 33).Value = 7.4
$ws.Cells.Item(4, 34).Value = 11.5
$ws.Cells.Item(4, 35).Value = 28
$ws.Cells.Item(4, 36).Value = 21
$ws.Cells.Item(4, 37).Value = 14
$ws.Cells.Item(4, 38).Value = 22
$ws.Cells.Item(4, 39).Value = 55
$ws.Cells.Item(4, 40).Value = 14
$ws.Cells.Item(4, 41).Value = 38

# Row 5: Portuguese Segunda Liga, Vizela vs Pacos Ferreira
$ws.Cells.Item(5, 1).Value = "Portuguese Segunda Liga"
$ws.Cells.Item(5, 3).Value = "17:15:00"
$ws.Cells.Item(5, 4).Value = "Vizela"
$ws.Cells.Item(5, 5).Value = "Pacos Ferreira"
$ws.Cells.Item(5, 6).Value = 1.78
$ws.Cells.Item(5, 7).Value = 1.8
$ws.Cells.Item(5, 8).Value = 5.5
$ws.Cells.Item(5, 9).Value = 5.8
$ws.Cells.Item(5, 10).Value = 3.75
$ws.Cells.Item(5, 11).Value = 3.9
$ws.Cells.Item(5, 12).Value = 1.46
$ws.Cells.Item(5, 13).Value = 1.08
$ws.Cells.Item(5, 14).Value = 3.35
$ws.Cells.Item(5, 15).Value = 1.41
$ws.Cells.Item(5, 16).Value = 1.79
$ws.Cells.Item(5, 17).Value = 2.2
$ws.Cells.Item(5, 18).Value = 1.29
$ws.Cells.Item(5, 19).Value = 4.1
$ws.Cells.Item(5, 20).Value = 2
$ws.Cells.Item(5, 21).Value = 1.89
$ws.Cells.Item(5, 22).Value = 1.21
$ws.Cells.Item(5, 23).Value = 2.24
$ws.Cells.Item(5, 24).Value = 12.5
$ws.Cells.Item(5, 25).Value = 16.5
$ws.Cells.Item(5, 26).Value = 48
$ws.Cells.Item(5, 27).Value = 160
$ws.Cells.Item(5, 28).Value = 7.4
$ws.Cells.Item(5, 29).Value = 8.4
$ws.Cells.Item(5, 30).Value = 22
$ws.Cells.Item(5, 31).Value = 110
$ws.Cells.Item(5, 32).Value = 9.6
$ws.Cells.Item(5, 33).Value = 9.4
$ws.Cells.Item(5, 34).Value = 24
$ws.Cells.Item(5, 35).Value = 260
$ws.Cells.Item(5, 36).Value = 19.5
$ws.Cells.Item(5, 37).Value = 21
$ws.Cells.Item(5, 38).Value = 44
$ws.Cells.Item(5, 39).Value = 580
$ws.Cells.Item(5, 40).Value = 16.5
$ws.Cells.Item(5, 41).Value = 160

# The old row 6 (Portuguese Segunda Liga / Vizela vs Pacos Ferreira placeholder) is gone;
# data now fits in rows 2-5, so remove the now-unused row 6.
$ws.Rows.Item(6).Delete()
